$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value2 = 111473774
$ws.Range("B3").Value2 = 89405
$ws.Range("D3").Value2 = 'NT'
$ws.Range("E3").Value2 = 1202
$ws.Range("F3").Value2 = 'Ullticka'
$ws.Range("G3").Value2 = 'Phellinidium ferrugineofuscum'
$ws.Range("H3").Value2 = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q3").Value2 = 703999.5190368021
$ws.Range("R3").Value2 = 6572850.823973293
$ws.Range("AO3").Value2 = 'granlåga'
# Row 4
$ws.Range("A4").Value2 = 111473777
$ws.Range("B4").Value2 = 89425
$ws.Range("D4").Value2 = 'NT'
$ws.Range("E4").Value2 = 5442
$ws.Range("F4").Value2 = 'Tallticka'
$ws.Range("G4").Value2 = 'Porodaedalea pini'
$ws.Range("H4").Value2 = '(Brot.) Murrill'
$ws.Range("Q4").Value2 = 704301.1177162804
$ws.Range("R4").Value2 = 6573209.392206083
$ws.Range("AO4").Value2 = 'gammeltall'
# Row 5
$ws.Range("A5").Value2 = 111473779
$ws.Range("B5").Value2 = 89425
$ws.Range("D5").Value2 = 'NT'
$ws.Range("E5").Value2 = 5442
$ws.Range("F5").Value2 = 'Tallticka'
$ws.Range("G5").Value2 = 'Porodaedalea pini'
$ws.Range("H5").Value2 = '(Brot.) Murrill'
$ws.Range("Q5").Value2 = 704193.4830821306
$ws.Range("R5").Value2 = 6572948.378178579
$ws.Range("AO5").Value2 = 'gammeltall'
# Row 6
$ws.Range("A6").Value2 = 111473793
$ws.Range("B6").Value2 = 93388
$ws.Range("D6").Value2 = 'LC'
$ws.Range("E6").Value2 = 2180
$ws.Range("F6").Value2 = 'Blåmossa'
$ws.Range("G6").Value2 = 'Leucobryum glaucum'
$ws.Range("H6").Value2 = '(Hedw.) Ångstr.'
$ws.Range("Q6").Value2 = 703959.3331032015
$ws.Range("R6").Value2 = 6572805.612961343
$ws.Range("AO6").ClearContents()
# Row 7
$ws.Range("A7").Value2 = 111473775
$ws.Range("Q7").Value2 = 703969.3444121893
$ws.Range("R7").Value2 = 6572791.287347207
# Row 8
$ws.Range("A8").Value2 = 111473783
$ws.Range("B8").Value2 = 89686
$ws.Range("E8").Value2 = 658
$ws.Range("F8").Value2 = 'Rosenticka'
$ws.Range("G8").Value2 = 'Rhodofomes roseus'
$ws.Range("H8").Value2 = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q8").Value2 = 703998.3853129407
$ws.Range("R8").Value2 = 6572852.813158008
# Row 9
$ws.Range("A9").Value2 = 111473776
$ws.Range("B9").Value2 = 89405
$ws.Range("E9").Value2 = 1202
$ws.Range("F9").Value2 = 'Ullticka'
$ws.Range("G9").Value2 = 'Phellinidium ferrugineofuscum'
$ws.Range("H9").Value2 = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q9").Value2 = 703970.8884549731
$ws.Range("R9").Value2 = 6572810.333898042
$ws.Range("AO9").Value2 = 'granlåga'
# Row 11
$ws.Range("A11").Value2 = 111473784
$ws.Range("B11").Value2 = 73634
$ws.Range("D11").Value2 = 'LC'
$ws.Range("E11").Value2 = 6426
$ws.Range("F11").Value2 = 'Kattfotslav'
$ws.Range("G11").Value2 = 'Felipes leucopellaeus'
$ws.Range("H11").Value2 = '(Ach.) Frisch & G.Thor'
$ws.Range("Q11").Value2 = 704135.470341172
$ws.Range("R11").Value2 = 6572843.267234835
$ws.Range("AO11").Value2 = 'äldre gran'
# Row 12
$ws.Range("A12").Value2 = 111473773
$ws.Range("B12").Value2 = 89405
$ws.Range("D12").Value2 = 'NT'
$ws.Range("E12").Value2 = 1202
$ws.Range("F12").Value2 = 'Ullticka'
$ws.Range("G12").Value2 = 'Phellinidium ferrugineofuscum'
$ws.Range("H12").Value2 = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("M12").ClearContents()
$ws.Range("Q12").Value2 = 704016.0051346947
$ws.Range("R12").Value2 = 6572801.994589122
$ws.Range("AO12").Value2 = 'granlåga'
# Row 13
$ws.Range("A13").Value2 = 111473792
$ws.Range("B13").Value2 = 5113
$ws.Range("D13").Value2 = 'LC'
$ws.Range("E13").Value2 = 100526
$ws.Range("F13").Value2 = 'Bronshjon'
$ws.Range("G13").Value2 = 'Callidium coriaceum'
$ws.Range("H13").Value2 = 'Paykull, 1800'
$ws.Range("M13").Value2 = 'äldre gnagspår'
$ws.Range("Q13").Value2 = 703965.55072247
$ws.Range("R13").Value2 = 6572785.445717536
$ws.Range("AO13").Value2 = 'torrgran'
# Row 14
$ws.Range("A14").Value2 = 111473791
$ws.Range("B14").Value2 = 93289
$ws.Range("D14").Value2 = 'LC'
$ws.Range("E14").Value2 = 2170
$ws.Range("F14").Value2 = 'Flagellkvastmossa'
$ws.Range("G14").Value2 = 'Dicranum flagellare'
$ws.Range("H14").Value2 = 'Hedw.'
$ws.Range("Q14").Value2 = 704004.9502936595
$ws.Range("R14").Value2 = 6572835.740028554
$ws.Range("AO14").Value2 = 'låga av tall'
